$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 2.458217522889299
$ws.Cells.Item(2, 4).Value = 0.2838200445432904
$ws.Cells.Item(2, 5).Value = 0.1396907984917775
$ws.Cells.Item(2, 6).Value = 6.80759003412939
$ws.Cells.Item(2, 7).Value = 0.002697211286910887
$ws.Cells.Item(2, 10).Value = 0.1417496444457811
$ws.Cells.Item(2, 12).Value = 1.038824214430605
$ws.Cells.Item(2, 13).Value = 0.7957860429121553
$ws.Cells.Item(3, 2).Value = 2.413157449983714
$ws.Cells.Item(3, 4).Value = 0.2544742934871635
$ws.Cells.Item(3, 5).Value = 0.1210603239397514
$ws.Cells.Item(3, 6).Value = 6.685570878255021
$ws.Cells.Item(3, 7).Value = 0.002707607811294551
$ws.Cells.Item(3, 10).Value = 0.1252084768248949
$ws.Cells.Item(3, 12).Value = 1.005668542131559
$ws.Cells.Item(3, 13).Value = 0.7768370955421844
$ws.Cells.Item(4, 2).Value = 2.387138399535701
$ws.Cells.Item(4, 4).Value = 0.2366359290947457
$ws.Cells.Item(4, 5).Value = 0.1096626335943824
$ws.Cells.Item(4, 6).Value = 6.614285832956796
$ws.Cells.Item(4, 7).Value = 0.002714314004745793
$ws.Cells.Item(4, 10).Value = 0.1150261744335381
$ws.Cells.Item(4, 12).Value = 0.9861727002578107
$ws.Cells.Item(4, 13).Value = 0.7657938199764871
$ws.Cells.Item(5, 2).Value = 2.37694849909451
$ws.Cells.Item(5, 4).Value = 0.2294097729609064
$ws.Cells.Item(5, 5).Value = 0.1050273848789232
$ws.Cells.Item(5, 6).Value = 6.586137278326277
$ws.Cells.Item(5, 7).Value = 0.002717128341331643
$ws.Cells.Item(5, 10).Value = 0.1108695739951315
$ws.Cells.Item(5, 12).Value = 0.9784431497423896
$ws.Cells.Item(5, 13).Value = 0.7614415989522385
$ws.Cells.Item(6, 2).Value = 2.375281385125646
$ws.Cells.Item(6, 4).Value = 0.2282124175077911
$ws.Cells.Item(6, 5).Value = 0.1042582435596771
$ws.Cells.Item(6, 6).Value = 6.581517242551797
$ws.Cells.Item(6, 7).Value = 0.002717600593240903
$ws.Cells.Item(6, 10).Value = 0.110178916935908
$ws.Cells.Item(6, 12).Value = 0.9771726147203594
$ws.Cells.Item(6, 13).Value = 0.7607278345765636
$ws.Cells.Item(7, 2).Value = 2.386999304146968
$ws.Cells.Item(7, 4).Value = 0.2365383028445649
$ws.Cells.Item(7, 5).Value = 0.1096000842851197
$ws.Cells.Item(7, 6).Value = 6.6139025824599
$ws.Cells.Item(7, 7).Value = 0.002714351629243719
$ws.Cells.Item(7, 10).Value = 0.1149701472509577
$ws.Cells.Item(7, 12).Value = 0.9860675876591927
$ws.Cells.Item(7, 13).Value = 0.7657345260669643
$ws.Cells.Item(8, 2).Value = 2.442337916623273
$ws.Cells.Item(8, 4).Value = 0.2736629469555112
$ws.Cells.Item(8, 5).Value = 0.1332577647106419
$ws.Cells.Item(8, 6).Value = 6.764755157231292
$ws.Cells.Item(8, 7).Value = 0.002700729256500365
$ws.Cells.Item(8, 10).Value = 0.1360512036086305
$ws.Cells.Item(8, 12).Value = 1.027212283025221
$ws.Cells.Item(8, 13).Value = 0.7891292356054933
$ws.Cells.Item(9, 2).Value = 2.564006322020532
$ws.Cells.Item(9, 4).Value = 0.3479947177460474
$ws.Cells.Item(9, 5).Value = 0.1800293600778531
$ws.Cells.Item(9, 6).Value = 7.090052917855331
$ws.Cells.Item(9, 7).Value = 0.002676559401487406
$ws.Cells.Item(9, 10).Value = 0.1772196429262465
$ws.Cells.Item(9, 12).Value = 1.11481487269063
$ws.Cells.Item(9, 13).Value = 0.8397370148142471
$ws.Cells.Item(10, 2).Value = 2.661533153017899
$ws.Cells.Item(10, 4).Value = 0.4036874852777999
$ws.Cells.Item(10, 5).Value = 0.2146969502910139
$ws.Cells.Item(10, 6).Value = 7.347934887929057
$ws.Cells.Item(10, 7).Value = 0.002660328820445812
$ws.Cells.Item(10, 10).Value = 0.2074135894358307
$ws.Cells.Item(10, 12).Value = 1.183516908222117
$ws.Cells.Item(10, 13).Value = 0.8798628462298268
$ws.Cells.Item(11, 2).Value = 2.707696293629738
$ws.Cells.Item(11, 4).Value = 0.4292915374342101
$ws.Cells.Item(11, 5).Value = 0.2305508958695981
$ws.Cells.Item(11, 6).Value = 7.469557505842829
$ws.Cells.Item(11, 7).Value = 0.002653271577311064
$ws.Cells.Item(11, 10).Value = 0.2211501418560715
$ws.Cells.Item(11, 12).Value = 1.215742393312979
$ws.Cells.Item(11, 13).Value = 0.8987706288010457
$ws.Cells.Item(12, 2).Value = 2.725437792860248
$ws.Cells.Item(12, 4).Value = 0.4390285094941362
$ws.Cells.Item(12, 5).Value = 0.236567718935575
$ws.Cells.Item(12, 6).Value = 7.516249295143723
$ws.Cells.Item(12, 7).Value = 0.002650645679166695
$ws.Cells.Item(12, 10).Value = 0.2263529479037345
$ws.Cells.Item(12, 12).Value = 1.228087492716583
$ws.Cells.Item(12, 13).Value = 0.9060257213727994
$ws.Cells.Item(13, 2).Value = 2.721605223645099
$ws.Cells.Item(13, 4).Value = 0.4369296010475807
$ws.Cells.Item(13, 5).Value = 0.2352712771916714
$ws.Cells.Item(13, 6).Value = 7.506164839500173
$ws.Cells.Item(13, 7).Value = 0.002651209150001645
$ws.Cells.Item(13, 10).Value = 0.2252323681033062
$ws.Cells.Item(13, 12).Value = 1.225422403559605
$ws.Cells.Item(13, 13).Value = 0.9044589596590242
$ws.Cells.Item(14, 2).Value = 2.709150664328831
$ws.Cells.Item(14, 4).Value = 0.4300917618381845
$ws.Cells.Item(14, 5).Value = 0.2310456299768902
$ws.Cells.Item(14, 6).Value = 7.473386028918355
$ws.Cells.Item(14, 7).Value = 0.002653054612623686
$ws.Cells.Item(14, 10).Value = 0.2215781530596956
$ws.Cells.Item(14, 12).Value = 1.216755173361605
$ws.Cells.Item(14, 13).Value = 0.8993655970823511
$ws.Cells.Item(15, 2).Value = 2.701555876111627
$ws.Cells.Item(15, 4).Value = 0.4259088403364331
$ws.Cells.Item(15, 5).Value = 0.2284590666969279
$ws.Cells.Item(15, 6).Value = 7.453391362406762
$ws.Cells.Item(15, 7).Value = 0.002654191060050836
$ws.Cells.Item(15, 10).Value = 0.219340007819028
$ws.Cells.Item(15, 12).Value = 1.211464803683157
$ws.Cells.Item(15, 13).Value = 0.8962581851045428
$ws.Cells.Item(16, 2).Value = 2.658552627054007
$ws.Cells.Item(16, 4).Value = 0.4020198342927586
$ws.Cells.Item(16, 5).Value = 0.213662645492434
$ws.Cells.Item(16, 6).Value = 7.340074672268059
$ws.Cells.Item(16, 7).Value = 0.002660796558823865
$ws.Cells.Item(16, 10).Value = 0.2065159734326301
$ws.Cells.Item(16, 12).Value = 1.181430649550151
$ws.Cells.Item(16, 13).Value = 0.8786404322714105
$ws.Cells.Item(17, 2).Value = 2.632633357634859
$ws.Cells.Item(17, 4).Value = 0.3874354439655292
$ws.Cells.Item(17, 5).Value = 0.2046077959346917
$ws.Cells.Item(17, 6).Value = 7.271673961944884
$ws.Cells.Item(17, 7).Value = 0.00266493208509228
$ws.Cells.Item(17, 10).Value = 0.1986497757422825
$ws.Cells.Item(17, 12).Value = 1.163256271119138
$ws.Cells.Item(17, 13).Value = 0.868000822481946
$ws.Cells.Item(18, 2).Value = 2.617894360034484
$ws.Cells.Item(18, 4).Value = 0.379072062713675
$ws.Cells.Item(18, 5).Value = 0.1994074654054003
$ws.Cells.Item(18, 6).Value = 7.232736454794718
$ws.Cells.Item(18, 7).Value = 0.002667341449717287
$ws.Cells.Item(18, 10).Value = 0.1941254012872093
$ws.Cells.Item(18, 12).Value = 1.152894259098076
$ws.Cells.Item(18, 13).Value = 0.8619427251882144
$ws.Cells.Item(19, 2).Value = 2.612932968851169
$ws.Cells.Item(19, 4).Value = 0.376244611942127
$ws.Cells.Item(19, 5).Value = 0.1976480208458753
$ws.Cells.Item(19, 6).Value = 7.219621920759977
$ws.Cells.Item(19, 7).Value = 0.002668162506342791
$ws.Cells.Item(19, 10).Value = 0.1925935082600887
$ws.Cells.Item(19, 12).Value = 1.149401495598596
$ws.Cells.Item(19, 13).Value = 0.8599020938409154
$ws.Cells.Item(20, 2).Value = 2.635374994680262
$ws.Cells.Item(20, 4).Value = 0.3889853533400753
$ws.Cells.Item(20, 5).Value = 0.2055708858909071
$ws.Cells.Item(20, 6).Value = 7.278913326843906
$ws.Cells.Item(20, 7).Value = 0.002664488674710799
$ws.Cells.Item(20, 10).Value = 0.1994871338573034
$ws.Cells.Item(20, 12).Value = 1.165181488894234
$ws.Cells.Item(20, 13).Value = 0.8691270522037087
$ws.Cells.Item(21, 2).Value = 2.712801783273164
$ws.Cells.Item(21, 4).Value = 0.4320990600007519
$ws.Cells.Item(21, 5).Value = 0.2322864342919928
$ws.Cells.Item(21, 6).Value = 7.482996571019839
$ws.Cells.Item(21, 7).Value = 0.002652511295849277
$ws.Cells.Item(21, 10).Value = 0.2226514481374835
$ws.Cells.Item(21, 12).Value = 1.219297077293589
$ws.Cells.Item(21, 13).Value = 0.9008590509723717
$ws.Cells.Item(22, 2).Value = 2.764924187522297
$ws.Cells.Item(22, 4).Value = 0.460518172389186
$ws.Cells.Item(22, 5).Value = 0.249824603415064
$ws.Cells.Item(22, 6).Value = 7.620091255161469
$ws.Cells.Item(22, 7).Value = 0.00264495440696825
$ws.Cells.Item(22, 10).Value = 0.2377973585199129
$ws.Cells.Item(22, 12).Value = 1.255493539792923
$ws.Cells.Item(22, 13).Value = 0.9221527272167975
$ws.Cells.Item(23, 2).Value = 2.736965749024648
$ws.Cells.Item(23, 4).Value = 0.4453273642943714
$ws.Cells.Item(23, 5).Value = 0.2404565830556891
$ws.Cells.Item(23, 6).Value = 7.546575963609484
$ws.Cells.Item(23, 7).Value = 0.002648962985253966
$ws.Cells.Item(23, 10).Value = 0.2297127810858939
$ws.Cells.Item(23, 12).Value = 1.236098235322089
$ws.Cells.Item(23, 13).Value = 0.9107367563436384
$ws.Cells.Item(24, 2).Value = 2.634134995083457
$ws.Cells.Item(24, 4).Value = 0.3882845729147562
$ws.Cells.Item(24, 5).Value = 0.2051354561256318
$ws.Cells.Item(24, 6).Value = 7.275639205850723
$ws.Cells.Item(24, 7).Value = 0.002664689041581573
$ws.Cells.Item(24, 10).Value = 0.1991085705277271
$ws.Cells.Item(24, 12).Value = 1.16431082810962
$ws.Cells.Item(24, 13).Value = 0.8686177008117397
$ws.Cells.Item(25, 2).Value = 2.529671927051311
$ws.Cells.Item(25, 4).Value = 0.3277070429963658
$ws.Cells.Item(25, 5).Value = 0.1673283674118835
$ws.Cells.Item(25, 6).Value = 6.998802069240128
$ws.Cells.Item(25, 7).Value = 0.002682828116338637
$ws.Cells.Item(25, 10).Value = 0.1660955831146822
$ws.Cells.Item(25, 12).Value = 1.090363033537869
$ws.Cells.Item(25, 13).Value = 0.8255337392173772
